$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = "bne0utiuzq@gmail.com"
$ws.Range("L3").Value = "bne0utiuzq@gmail.com"
